$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 19519.715
$ws.Range("I21").Value = 22333.334
$ws.Range("J21").Value = 17409.5
$ws.Range("K21").Value = 22333.334
$ws.Range("L21").Value = 17409.5
$ws.Range("M21").Value = -21865.334
$ws.Range("N21").Value = -18345.5
$ws.Range("H23").Value = 19519.715
$ws.Range("I23").Value = 22333.334
$ws.Range("J23").Value = 17409.5
$ws.Range("K23").Value = 22333.334
$ws.Range("L23").Value = 17409.5
$ws.Range("M23").Value = -22099.334
$ws.Range("N23").Value = -17877.5
$ws.Range("H41").Value = 1259.55
$ws.Range("I41").Value = 1033.4166
$ws.Range("J41").Value = 1598.75
$ws.Range("K41").Value = 1033.4166
$ws.Range("L41").Value = 1598.75
$ws.Range("M41").Value = -593.4166
$ws.Range("N41").Value = -2478.75
$ws.Range("H54").Value = 3895.1428
$ws.Range("J54").Value = 5000
$ws.Range("L54").Value = 5000
$ws.Range("N54").Value = -5972
$ws.Range("H64").Value = 4149.8
$ws.Range("I64").Value = 4263
$ws.Range("K64").Value = 4263
$ws.Range("M64").Value = -4015
$ws.Range("H67").Value = 4149.8
$ws.Range("I67").Value = 4263
$ws.Range("K67").Value = 4263
$ws.Range("M67").Value = -3405
$ws.Range("H70").Value = 989
$ws.Range("I70").Value = 999
$ws.Range("K70").Value = 2997
$ws.Range("M70").Value = -2727
$ws.Range("H73").Value = 989
$ws.Range("I73").Value = 999
$ws.Range("K73").Value = 2997
$ws.Range("M73").Value = -2061
$ws.Range("H86").Value = 3268.4
$ws.Range("I86").Value = 3606.8
$ws.Range("J86").Value = 2930
$ws.Range("K86").Value = 3606.8
$ws.Range("L86").Value = 2930
$ws.Range("M86").Value = -2483.8
$ws.Range("N86").Value = -5176
$ws.Range("H89").Value = 3268.4
$ws.Range("I89").Value = 3606.8
$ws.Range("J89").Value = 2930
$ws.Range("K89").Value = 18034
$ws.Range("L89").Value = 14650
$ws.Range("M89").Value = -12418
$ws.Range("N89").Value = -25882
$ws.Range("H137").Value = 1452.9803
$ws.Range("I137").Value = 1167.875
$ws.Range("J137").Value = 1933.1578
$ws.Range("K137").Value = 3503.625
$ws.Range("L137").Value = 5799.4734
$ws.Range("M137").Value = -953.625
$ws.Range("N137").Value = -10899.4734
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2306.4167
$ws.Range("I132").Value = 1597.7142
$ws.Range("J132").Value = 3960.0557
$ws.Range("K132").Value = 4793.142599999999
$ws.Range("L132").Value = 11880.1671
$ws.Range("M132").Value = -2263.142599999999
$ws.Range("N132").Value = -16940.1671
$ws.Range("H134").Value = 34080
$ws.Range("J134").Value = 34080
$ws.Range("L134").Value = 34080
$ws.Range("N134").Value = -44220
$ws.Range("H141").Value = 34868.5
$ws.Range("J141").Value = 34868.5
$ws.Range("L141").Value = 34868.5
$ws.Range("N141").Value = -45228.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 19954.75
$ws.Range("J58").Value = 23273
$ws.Range("L58").Value = 23273
$ws.Range("N58").Value = -23861
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 11769318
$ws.Range("I62").Value = 6440
$ws.Range("J62").Value = 28573428
$ws.Range("K62").Value = 6440
$ws.Range("L62").Value = 28573428
$ws.Range("M62").Value = -5816
$ws.Range("N62").Value = -28574676
$ws.Range("H65").Value = 11769318
$ws.Range("I65").Value = 6440
$ws.Range("J65").Value = 28573428
$ws.Range("K65").Value = 32200
$ws.Range("L65").Value = 142867140
$ws.Range("M65").Value = -29080
$ws.Range("N65").Value = -142873380
$ws.Range("H134").Value = 16668178
$ws.Range("I134").Value = 1570.6818
$ws.Range("K134").Value = 4712.0454
$ws.Range("M134").Value = -2177.0454
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 10843.1875
$ws.Range("J63").Value = 17064.334
$ws.Range("L63").Value = 51193.00199999999
$ws.Range("N63").Value = -52691.00199999999
$ws.Range("H66").Value = 10843.1875
$ws.Range("J66").Value = 17064.334
$ws.Range("L66").Value = 153579.006
$ws.Range("N66").Value = -161067.006
$ws.Range("H106").Value = 3241.8572
$ws.Range("J106").Value = 3276.3333
$ws.Range("L106").Value = 9828.999899999999
$ws.Range("N106").Value = -11720.9999
$ws.Range("H131").Value = 22762120
$ws.Range("J131").Value = 47735.22
$ws.Range("L131").Value = 143205.66
$ws.Range("N131").Value = -153285.66
$ws.Range("H140").Value = 2872.5957
$ws.Range("I140").Value = 2002.579
$ws.Range("J140").Value = 3462.9644
$ws.Range("K140").Value = 6007.737
$ws.Range("L140").Value = 10388.8932
$ws.Range("M140").Value = -827.7370000000001
$ws.Range("N140").Value = -20748.8932
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 773.8333
$ws.Range("J22").Value = 699.1667
$ws.Range("L22").Value = 699.1667
$ws.Range("N22").Value = -1289.1667
$ws.Range("H27").Value = 773.8333
$ws.Range("J27").Value = 699.1667
$ws.Range("L27").Value = 699.1667
$ws.Range("N27").Value = -913.1667
$ws.Range("H46").Value = 5700
$ws.Range("J46").Value = 6875
$ws.Range("L46").Value = 6875
$ws.Range("N46").Value = -7251
$ws.Range("H55").Value = 283.65
$ws.Range("J55").Value = 458.33334
$ws.Range("L55").Value = 458.33334
$ws.Range("N55").Value = -804.33334
$ws.Range("I122").Value = 83334670
$ws.Range("K122").Value = 250004010
$ws.Range("M122").Value = -250001560
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 100002400
$ws.Range("I62").Value = 125002250
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 125002250
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -125001626
$ws.Range("N62").Value = -4248
$ws.Range("H65").Value = 100002400
$ws.Range("I65").Value = 125002250
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 625011250
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -625008130
$ws.Range("N65").Value = -21240
$ws.Range("H136").Value = 1003.3214
$ws.Range("I136").Value = 716.86957
$ws.Range("K136").Value = 2150.60871
$ws.Range("M136").Value = 399.39129
